$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = -20.6
$ws.Range("A7").Value = -21.286
$ws.Range("B7").Value = 6.312
$ws.Range("B12").Value = 5.401999999999999
$ws.Range("B15").Value = 5.359000000000001
$ws.Range("A16").Value = -21.539
$ws.Range("A20").Value = -21.795
$ws.Range("B20").Value = 5.619000000000001
$ws.Range("B21").Value = 8.904
$ws.Range("B22").Value = 7.034000000000001
$ws.Range("B23").Value = 7.359999999999999
$ws.Range("A28").Value = -21.878
$ws.Range("A29").Value = -21.675
$ws.Range("B29").Value = 5.731
$ws.Range("A32").Value = -21.705
$ws.Range("B34").Value = 8.059000000000001
$ws.Range("A40").Value = -20.312
$ws.Range("B42").Value = 7.238000000000001
$ws.Range("B43").Value = 5.529000000000001
$ws.Range("B44").Value = 5.08
$ws.Range("B45").Value = 5.315
$ws.Range("A46").Value = -20.832
$ws.Range("B46").Value = 7.604000000000001
$ws.Range("B50").Value = 5.528
$ws.Range("A51").Value = -20.771
$ws.Range("B51").Value = 7.779000000000001
$ws.Range("A52").Value = -21.316
$ws.Range("A57").Value = -22.263
$ws.Range("A59").Value = -22.395
$ws.Range("A62").Value = -21.785
$ws.Range("A66").Value = -21.504
$ws.Range("B66").Value = 5.473
$ws.Range("B67").Value = 5.194999999999999
$ws.Range("A73").Value = -20.53
$ws.Range("A74").Value = -21.043
$ws.Range("B79").Value = 5.681
$ws.Range("B84").Value = 5.781000000000001
$ws.Range("A92").Value = -21.333
$ws.Range("B92").Value = 5.526999999999999
$ws.Range("B97").Value = 6.296000000000001
$ws.Range("A100").Value = -21.481
